$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Columns B:C had only a (visually no-op) column-wide style and an
#     explicit width equal to the sheet default - drop that declaration
#     entirely, then restore the few per-cell formats in B:C that really
#     do matter (header/footer borders + the bold sub-header row).
$ws.Columns("B:C").ClearFormats()

$ws.Range("B1:C1").Borders.Item(9).LineStyle = 1
$ws.Range("B1:C1").Borders.Item(9).Weight = 2

$ws.Range("B2:C2").Font.Bold = $true

$ws.Range("B42:C42").Borders.Item(9).LineStyle = 1
$ws.Range("B42:C42").Borders.Item(9).Weight = -4138

# --- 2. Drop the same no-op border style from the rest of the data area
#     (columns A, D, E - B/C already handled above).
$ws.Range("A3:A41").Style = "Normal"
$ws.Range("D3:E41").Style = "Normal"
$ws.Range("C44").Style = "Normal"

# Re-apply wrap-text to the column-A question cells that need it (the
# blanket style reset above cleared it).
$wrapRows = @(5,6,8,12,14,15,17)
foreach ($r in $wrapRows) {
    $ws.Range("A$r").WrapText = $true
}

# E42 also loses its border entirely (unlike the rest of row 42).
$ws.Range("E42").Style = "Normal"

# --- 3. Flip the response-option coding for the reverse-scored SCQ items:
#     "0= No, 1= Yes"  ->  "1=No, 0=Yes"
$invertedRows = @(4,11,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42)
foreach ($r in $invertedRows) {
    $ws.Range("E$r").Value = "1=No, 0=Yes"
}

# --- 4. Update the window/view state to match the new edit position.
$ws.Range("A31").Select()
$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("A15")
